# Auto-generated script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.709.57'
$ws.Range('E2').Value = '  +2.60%  '
$ws.Range('D3').Value = '3.008.33'
$ws.Range('E3').Value = '  +2.31%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '511.67'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +6.29%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.22'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +7.49%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.435'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +4.79%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '7.56'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +8.76%  '
$ws.Range('E10').Value = '  +9.66%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.358'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +3.79%  '
$ws.Range('E12').Value = '  +3.05%  '
$ws.Range('D13').Value = '3.522.76'
$ws.Range('E13').Value = '  +2.31%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '25.82'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +7.12%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000158'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +13.94%  '
$ws.Range('D16').Value = '56.764.00'
$ws.Range('E16').Value = '  +2.94%  '
$ws.Range('D17').Value = '3.006.48'
$ws.Range('E17').Value = '  +2.54%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '5.94'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +7.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.56'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +5.80%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.88'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +7.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '331.93'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +7.51%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.484'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +6.04%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.04'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.172'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +10.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '0.0₃0909'
$ws.Range('E27').Value = '  +8.46%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.73'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +3.86%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.10'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +10.72%  '
$ws.Range('E30').Value = '  +7.42%  '
$ws.Range('E31').Value = '  +8.19%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '20.74'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +8.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '154.38'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +5.46%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.60'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +6.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.70'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.28'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.27%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0680'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +6.37%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.85'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').Value = '3.039.82'
$ws.Range('E39').Value = '  +2.37%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '36.98'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +3.43%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('E42').Value = '  +4.29%  '
$ws.Range('D43').Value = '2.276.00'
$ws.Range('E43').Value = '  +8.09%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.69'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +5.96%  '
$ws.Range('B45').Value = 'ONDO'
$ws.Range('C45').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.01'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.97%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.42'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +4.94%  '
$ws.Range('E47').Value = '  +18.86%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0241'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +6.30%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '5.85'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +6.40%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.66'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +5.23%  '
$ws.Range('E51').Value = '  +7.25%  '

Write-Host "Applied 94 cell updates (32 numeric-text, 62 plain text)"
